$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.027.66'
$ws.Range("E2").Value = '  +3.16%  '
$ws.Range("D3").Value = '1.688.37'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.12'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.264'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.75%  '
$ws.Range("E10").Value = '  +0.20%  '
$ws.Range("E11").Value = '  -0.51%  '
$ws.Range("D12").Value = '1.929.05'
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("D13").Value = '1.698.10'
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.88'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '250.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.09%  '
$ws.Range("D18").Value = '28.006.25'
$ws.Range("E18").Value = '  +3.03%  '
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("E20").Value = '  -3.44%  '
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.66%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.33'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.48'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.46%  '
$ws.Range("E28").Value = '  +0.26%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.26'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.30%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("E32").Value = '  +0.21%  '
$ws.Range("E33").Value = '  -1.90%  '
$ws.Range("D34").Value = '1.440.80'
$ws.Range("E34").Value = '  -6.97%  '
$ws.Range("E35").Value = '  -2.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.949'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("E37").Value = '  +0.19%  '
$ws.Range("E38").Value = '  -2.03%  '
$ws.Range("E39").Value = '  +0.20%  '
$ws.Range("E40").Value = '  -2.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.52'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.69%  '
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D44").Value = '1.836.15'
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("E45").Value = '  -1.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.796'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.61%  '
$ws.Range("E47").Value = '  +6.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '89.41'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.34%  '
$ws.Range("E49").Value = '  -1.00%  '
$ws.Range("E50").Value = '  -1.22%  '
$ws.Range("E51").Value = '  -4.13%  '
